$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (tab) - matches the <sheet name="..."> change
$ws.Name = "Through 2021-11-13"

# Row 13: "November (through 11-12)" -> "November (through 11-13)" and updated values
$ws.Range("A13").Value = "November (through 11-13)"
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 14
$ws.Range("D13").Value = 0.0667
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 29
$ws.Range("G13").Value = 0.0938
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = 58
$ws.Range("J13").Value = 0.0169
$ws.Range("K13").Value = 5
$ws.Range("L13").Value = 23
$ws.Range("M13").Value = 0.1786
$ws.Range("N13").Value = 4
$ws.Range("O13").Value = 18
$ws.Range("P13").Value = 0.1818
$ws.Range("Q13").Value = 2
$ws.Range("R13").Value = 79
$ws.Range("S13").Value = 0.0247
$ws.Range("T13").Value = 1
$ws.Range("U13").Value = 86
$ws.Range("V13").Value = 0.0115

# Row 14: Total row updated values
$ws.Range("B14").Value = 33
$ws.Range("C14").Value = 240
$ws.Range("D14").Value = 0.1209
$ws.Range("E14").Value = 55
$ws.Range("F14").Value = 463
$ws.Range("G14").Value = 0.1062
$ws.Range("H14").Value = 62
$ws.Range("I14").Value = 707
$ws.Range("J14").Value = 0.0806
$ws.Range("K14").Value = 71
$ws.Range("L14").Value = 572
$ws.Range("M14").Value = 0.1104
$ws.Range("N14").Value = 52
$ws.Range("O14").Value = 452
$ws.Range("P14").Value = 0.1032
$ws.Range("Q14").Value = 56
$ws.Range("R14").Value = 1082
$ws.Range("S14").Value = 0.0492
$ws.Range("T14").Value = 89
$ws.Range("U14").Value = 1440
$ws.Range("V14").Value = 0.0582
